# Generate Report for handoff
#
# The localization-status report is regenerated:
#  - The file "1247372a-9347-4ef6-bb28-6a6774aa4064.md" moved from
#    "Ready for handoff" to "In Translation" on every sheet (Overview,
#    zh-cn, de-de).
#  - The other rows that were still "Ready for handoff" got a freshly
#    stamped "Latest Handoff Datetime" (zh-cn -> 2016-01-26 12:35:54,
#    de-de -> 2016-01-26 12:36:07), except the de-de row for the file
#    above, whose handoff datetime was not refreshed (status changed
#    instead).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B4").Value = "In Translation"
$ov.Range("C4").Value = "In Translation"

# ---- zh-cn sheet -------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("B4").Value = "In Translation"
$zh.Range("D2").Value = "2016-01-26 12:35:54"
$zh.Range("D3").Value = "2016-01-26 12:35:54"
$zh.Range("D4").Value = "2016-01-26 12:35:54"
$zh.Range("D5").Value = "2016-01-26 12:35:54"
$zh.Range("D6").Value = "2016-01-26 12:35:54"
$zh.Range("D7").Value = "2016-01-26 12:35:54"

# ---- de-de sheet -------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
$de.Range("B4").Value = "In Translation"
$de.Range("D2").Value = "2016-01-26 12:36:07"
$de.Range("D3").Value = "2016-01-26 12:36:07"
$de.Range("D5").Value = "2016-01-26 12:36:07"
$de.Range("D6").Value = "2016-01-26 12:36:07"
$de.Range("D7").Value = "2016-01-26 12:36:07"
